$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = -12.294
$ws.Range("A3").Value = -21.779
$ws.Range("E3").Value = 16.492
$ws.Range("E12").Value = 17.652
$ws.Range("A14").Value = -21.863
$ws.Range("A16").Value = -21.945
$ws.Range("C18").Value = -11.54
$ws.Range("A21").Value = -20.066
$ws.Range("A23").Value = -20.198
$ws.Range("C24").Value = -12.591
$ws.Range("E24").Value = 16.954
$ws.Range("A25").Value = -21.646
$ws.Range("C25").Value = -11.969
$ws.Range("E25").Value = 17.053
$ws.Range("A26").Value = -21.114
$ws.Range("C27").Value = -13.309
$ws.Range("A29").Value = -21.219
$ws.Range("C30").Value = -13.133
$ws.Range("C31").Value = -13.223
$ws.Range("C39").Value = -12.728
$ws.Range("A40").Value = -20.197
$ws.Range("E41").Value = 16.554
$ws.Range("C42").Value = -12.852
$ws.Range("C48").Value = -11.422
$ws.Range("E50").Value = 16.305
$ws.Range("C51").Value = -11.133
$ws.Range("C52").Value = -11.601
$ws.Range("A53").Value = -21.845
$ws.Range("E53").Value = 16.781
$ws.Range("C55").Value = -13.251
$ws.Range("C56").Value = -13.445
$ws.Range("E56").Value = 16.378
$ws.Range("A57").Value = -22.1
$ws.Range("C57").Value = -13.813
$ws.Range("E57").Value = 16.408
$ws.Range("E58").Value = 16.545
$ws.Range("A59").Value = -22.323
$ws.Range("C60").Value = -13.048
$ws.Range("E61").Value = 16.82
$ws.Range("E63").Value = 17.652
$ws.Range("E64").Value = 17.39
$ws.Range("A65").Value = -21.421
$ws.Range("A69").Value = -21.783
$ws.Range("E70").Value = 17.64
$ws.Range("E72").Value = 17.046
$ws.Range("C73").Value = -12.601
$ws.Range("C74").Value = -12.462
$ws.Range("A83").Value = -21.938
$ws.Range("E86").Value = 16.428
$ws.Range("C89").Value = -10.676
$ws.Range("E89").Value = 17.14
$ws.Range("C90").Value = -12.957
$ws.Range("A91").Value = -21.527
$ws.Range("C92").Value = -11.085
$ws.Range("A93").Value = -21.356
$ws.Range("E98").Value = 16.086
$ws.Range("A100").Value = -22.016
$ws.Range("E100").Value = 16.393
$ws.Range("E102").Value = 16.513
